# Add a "2022-Q4" fund-holdings sheet right after the "总计" summary sheet,
# and insert a matching 2022-Q4 row at the top of the "总计" table (pushing
# every existing quarter down by one row / one sheet position).

$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) "总计" sheet: insert new row 2, bump the old index column (A) by 1,
#    then fill in the 2022-Q4 figures.
# ---------------------------------------------------------------------
$summary.Rows.Item(2).Insert()

# Row-insert copies formatting down from the header row for the new row;
# strip that back to the default (unstyled) look used by every other
# data row in B:D.
$summary.Range("B2:D2").ClearFormats()

# The old rows 2..8 are now rows 3..9; their "A" index values (0..6) need
# to become 1..7 since a new 0-th entry was inserted above them.
for ($r = 9; $r -ge 3; $r--) {
    $old = [double]$summary.Cells.Item($r, 1).Value2
    $summary.Cells.Item($r, 1).Value = $old + 1
}

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q4"
$summary.Cells.Item(2, 3).Value = 25
$summary.Cells.Item(2, 4).Value = 2.04

# A2 needs the same bold/centered/bordered look as the rest of column A
# (A3:A9); copy format only from A3 onto A2 (value is untouched by a
# formats-only paste).
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2) Add the new "2022-Q4" worksheet right after "总计".
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $summary)
$q4.Name = "2022-Q4"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $q4.Cells.Item(1, $c + 2).Value = $headers[$c]
}

# code, name, scale, stockPosition, positionRatio, marketValue, rank
$data = @(
    @("011069", "工银成长精选混合A", "12.99", "91.27", "3.04", "0.3949", 7),
    @("005228", "汇添富港股通专注成长混合", "7.51", "85.64", "4.72", "0.3545", 5),
    @("013123", "汇添富精选核心优势一年持有混合A", "5.58", "83.43", "4.24", "0.2366", 7),
    @("010701", "恒越内需驱动混合A", "6.59", "90.58", "3.22", "0.2122", 9),
    @("008227", "宝盈研究精选混合A", "3.93", "91.85", "4.75", "0.1867", 9),
    @("671010", "西部利得策略优选混合A", "1.88", "92.90", "7.15", "0.1344", 6),
    @("013550", "汇添富品牌价值一年持有混合A", "2.24", "75.70", "4.25", "0.0952", 7),
    @("010702", "恒越内需驱动混合C", "1.98", "90.58", "3.22", "0.0638", 9),
    @("011070", "工银成长精选混合C", "1.82", "91.27", "3.04", "0.0553", 7),
    @("008228", "宝盈研究精选混合C", "1.01", "91.85", "4.75", "0.0480", 9),
    @("013028", "恒越品质生活混合", "1.25", "90.58", "3.12", "0.0390", 6),
    @("005143", "中融沪港深大消费主题灵活配置混合C", "0.69", "77.26", "5.50", "0.0380", 4),
    @("011060", "西部利得策略优选混合C", "0.48", "92.90", "7.15", "0.0343", 6),
    @("014220", "恒越医疗健康精选混合A", "0.73", "90.68", "3.40", "0.0248", 7),
    @("006241", "中融医疗健康精选混合C", "0.55", "93.89", "4.24", "0.0233", 4),
    @("013441", "西藏东财创新医疗六个月定开混合", "0.57", "82.77", "3.78", "0.0215", 10),
    @("005142", "中融沪港深大消费主题灵活配置混合A", "0.31", "77.26", "5.50", "0.0170", 4),
    @("005701", "上投摩根香港精选港股通混合A", "0.48", "89.99", "3.17", "0.0152", 5),
    @("013551", "汇添富品牌价值一年持有混合C", "0.28", "75.70", "4.25", "0.0119", 7),
    @("013124", "汇添富精选核心优势一年持有混合C", "0.27", "83.43", "4.24", "0.0114", 7),
    @("014221", "恒越医疗健康精选混合C", "0.31", "90.68", "3.40", "0.0105", 7),
    @("009734", "创金合信港股通大消费精选股票C", "0.19", "81.61", "3.63", "0.0069", 5),
    @("006240", "中融医疗健康精选混合A", "0.10", "93.89", "4.24", "0.0042", 4),
    @("009733", "创金合信港股通大消费精选股票A", "0.09", "81.61", "3.63", "0.0033", 5),
    @("016921", "上投摩根香港精选港股通混合C", "0.02", "89.99", "3.17", "0.0006", 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2

    $q4.Cells.Item($r, 1).Value = $i

    # Fund code (leading zeros) and the scale/position/ratio/value columns
    # are numeric-looking strings that must stay text -- force a Text
    # number format before assignment, otherwise "011069" becomes 11069.
    $q4.Cells.Item($r, 2).NumberFormat = "@"
    $q4.Cells.Item($r, 2).Value = $row[0]

    $q4.Cells.Item($r, 3).Value = $row[1]

    $q4.Cells.Item($r, 4).NumberFormat = "@"
    $q4.Cells.Item($r, 4).Value = $row[2]
    $q4.Cells.Item($r, 5).NumberFormat = "@"
    $q4.Cells.Item($r, 5).Value = $row[3]
    $q4.Cells.Item($r, 6).NumberFormat = "@"
    $q4.Cells.Item($r, 6).Value = $row[4]
    $q4.Cells.Item($r, 7).NumberFormat = "@"
    $q4.Cells.Item($r, 7).Value = $row[5]

    $q4.Cells.Item($r, 8).Value = $row[6]
}

# The Text-number-format trick above leaves an explicit "@" style on
# B,D,E,F,G; a formats-only paste from an unstyled cell (H1 is still
# blank/default at this point) restores the default look without
# touching any values.
$summary.Range("B10").Copy()
$q4.Range("B2:B26").PasteSpecial(-4122)
$q4.Range("D2:G26").PasteSpecial(-4122)

# Header row (B1:H1) and index column (A2:A26) get the bold/centered/
# bordered style used throughout the workbook -- copy it (format only)
# from the equivalent cells on "总计".
$summary.Range("B1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)

$summary.Range("A2").Copy()
$q4.Range("A2:A26").PasteSpecial(-4122)

Write-Output "done"
